# Update cryptocurrency price/volume data (cryptos.xlsx)
# Commit: Updated cryptos list on Sun Sep 22 19:27:53 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = "'62.850.46"
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = "'2.568.27"
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'582.44"
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').Value = "'143.64"
$ws.Range('E6').Value = '  -2.61%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = "'0.588"
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('E9').Value = '  -2.58%  '
$ws.Range('D10').Value = "'5.57"
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').Value = "'0.348"
$ws.Range('D13').Value = "'26.94"
$ws.Range('E13').Value = '  -2.00%  '
$ws.Range('D14').Value = "'3.031.65"
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').Value = "'62.744.54"
$ws.Range('E15').Value = '  -0.73%  '
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('D17').Value = "'2.573.16"
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = "'11.04"
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').Value = "'340.09"
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('D20').Value = "'4.32"
$ws.Range('E20').Value = '  -2.10%  '
$ws.Range('D21').Value = "'6.61"
$ws.Range('E21').Value = '  -3.78%  '
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').Value = "'5.73"
$ws.Range('E23').Value = '  +3.20%  '
$ws.Range('D24').Value = "'67.69"
$ws.Range('E24').Value = '  +1.21%  '
$ws.Range('D25').Value = "'1.58"
$ws.Range('E25').Value = '  +6.44%  '
$ws.Range('D26').Value = "'1.58"
$ws.Range('E26').Value = '  -3.27%  '
$ws.Range('E27').Value = '  -3.60%  '
$ws.Range('D28').Value = "'7.96"
$ws.Range('E28').Value = '  -2.17%  '
$ws.Range('E29').Value = '  -1.52%  '
$ws.Range('D30').Value = "'8.21"
$ws.Range('E30').Value = '  -3.32%  '
$ws.Range('E31').Value = '  -2.54%  '
$ws.Range('D32').Value = "'458.86"
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('D33').Value = "'0.0₃0793"
$ws.Range('E33').Value = '  -3.88%  '
$ws.Range('D34').Value = "'1.65"
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('D35').Value = "'176.66"
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('E37').Value = '  -2.45%  '
$ws.Range('D38').Value = "'18.80"
$ws.Range('E38').Value = '  -2.25%  '
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('D41').Value = "'1.68"
$ws.Range('E41').Value = '  -3.44%  '
$ws.Range('D42').Value = "'39.96"
$ws.Range('E42').Value = '  +0.88%  '
$ws.Range('D43').Value = "'157.20"
$ws.Range('E43').Value = '  +3.96%  '
$ws.Range('E44').Value = '  -3.54%  '
$ws.Range('D45').Value = "'21.15"
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('D46').Value = "'0.631"
$ws.Range('E46').Value = '  +2.78%  '
$ws.Range('D47').Value = "'0.0533"
$ws.Range('E47').Value = '  -2.92%  '
$ws.Range('D48').Value = "'0.0957"
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('D49').Value = "'0.0234"
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('D50').Value = "'17.99"
$ws.Range('E50').Value = '  -2.69%  '
$ws.Range('E51').Value = '  +0.06%  '
